# Subject Information Form - rebuild sheet per target revision
# (new "Scheduled" / "Actural" scan date+time columns, new subject rows, updated headers)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Wipe existing content/formatting/merges so we can rebuild cleanly.
# ---------------------------------------------------------------------------
$ws.Cells.UnMerge()
$ws.Cells.Clear()

# ---------------------------------------------------------------------------
# 2. Column widths (character units) -- matches target <cols> widths as
#    closely as this engine's pixel/sixth-of-a-character quantization allows.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 9
$ws.Columns("B").ColumnWidth = 16
$ws.Columns("C").ColumnWidth = 9
$ws.Columns("D").ColumnWidth = 9
$ws.Columns("E").ColumnWidth = 39
$ws.Columns("F").ColumnWidth = 39
$ws.Columns("G").ColumnWidth = 25
$ws.Columns("H").ColumnWidth = 25
$ws.Columns("I").ColumnWidth = 24
$ws.Columns("J").ColumnWidth = 17
$ws.Columns("K").ColumnWidth = 17
$ws.Columns("L").ColumnWidth = 17
$ws.Columns("M").ColumnWidth = 20
$ws.Columns("N").ColumnWidth = 20
$ws.Columns("O").ColumnWidth = 11
$ws.Columns("P").ColumnWidth = 36.8
$ws.Columns("Q").ColumnWidth = 11

# ---------------------------------------------------------------------------
# 3. Base font/alignment for the whole working area (14pt, centered) --
#    matches the sheet's "font1" look used throughout.
# ---------------------------------------------------------------------------
$all = $ws.Range("A1:Q6")
$all.Font.Size = 14
$all.HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. Row 1 -- section banner cells, with vertical divider borders at the
#    boundaries of each merged block, then merge.
# ---------------------------------------------------------------------------
$ws.Range("F1").Borders.Item(10).LineStyle = 1
$ws.Range("G1").Borders.Item(7).LineStyle = 1
$ws.Range("L1").Borders.Item(10).LineStyle = 1
$ws.Range("P1").Borders.Item(10).LineStyle = 1

$ws.Range("A1:F1").Merge()
$ws.Range("G1:L1").Merge()
$ws.Range("M1:P1").Merge()

$ws.Range("A1").Value = "Subject Information"
$ws.Range("G1").Value = "Things to confirm before scan"
$ws.Range("M1").Value = "Things to confirm after scan"

# ---------------------------------------------------------------------------
# 5. Row 2 -- column headers.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "ID"
$ws.Range("B2").Value = "Name"
$ws.Range("C2").Value = "Age"
$ws.Range("D2").Value = "Gender"
$ws.Range("E2").Value = "Subject Email"
$ws.Range("F2").Value = "Parent Email"
$ws.Range("G2").Value = "Scheduled Scan Date"
$ws.Range("H2").Value = "Scheduled Starting Time"
$ws.Range("I2").Value = "Subject Assent/Consent"
$ws.Range("J2").Value = "Parent Consent"
$ws.Range("K2").Value = "Payment Receipt"
$ws.Range("L2").Value = "Questionnaires"
$ws.Range("M2").Value = "Actural Scan Date"
$ws.Range("N2").Value = "Actural Scan Time"
$ws.Range("O2").Value = "Daris ID"
$ws.Range("P2").Value = "Notes"

# Vertical divider borders at the end of each logical header group.
$ws.Range("F2").Borders.Item(10).LineStyle = 1
$ws.Range("L2").Borders.Item(10).LineStyle = 1
$ws.Range("P2").Borders.Item(10).LineStyle = 1

# ---------------------------------------------------------------------------
# 6. Row 3 -- subject 1 (Maddy Coates).
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Maddy Coates"
$ws.Range("C3").Value = 18
$ws.Range("D3").Value = "Female"
$ws.Range("E3").Value = "mkcoates@student.unimelb.edu.au"
$ws.Range("G3").Value = 45274
$ws.Range("G3").NumberFormat = "mm-dd-yy"
$ws.Range("H3").Value = 0.33333333333333331
$ws.Range("H3").NumberFormat = "h:mm"
$ws.Range("I3").Value = "/"
$ws.Range("J3").Value = "/"
$ws.Range("K3").Value = "/"
$ws.Range("L3").Value = "?"
$ws.Range("M3").Value = 45274
$ws.Range("M3").NumberFormat = "d-mmm-yy"
$ws.Range("N3").Value = "8:00-10:10"
$ws.Range("P3").Value = "Sound wasn't on for the safty tasks"

$ws.Range("F3").Borders.Item(10).LineStyle = 1
$ws.Range("L3").Borders.Item(10).LineStyle = 1
$ws.Range("P3").Borders.Item(10).LineStyle = 1

# ---------------------------------------------------------------------------
# 7. Row 4 -- subject 2 (Ryan White).
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Ryan White"
$ws.Range("C4").Value = 16
$ws.Range("D4").Value = "Male"
$ws.Range("F4").Value = "whiteantsrule@gmail.com"
$ws.Range("G4").Value = 45306
$ws.Range("G4").NumberFormat = "mm-dd-yy"
$ws.Range("H4").Value = 0.4375
$ws.Range("H4").NumberFormat = "h:mm"

$ws.Range("F4").Borders.Item(10).LineStyle = 1
$ws.Range("L4").Borders.Item(10).LineStyle = 1
$ws.Range("P4").Borders.Item(10).LineStyle = 1

# ---------------------------------------------------------------------------
# 8. Row 5 -- subject 3 (Ethan White).
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Ethan White"
$ws.Range("C5").Value = 16
$ws.Range("D5").Value = "Male"
$ws.Range("F5").Value = "whiteantsrule@gmail.com"
$ws.Range("G5").Value = 45306
$ws.Range("G5").NumberFormat = "mm-dd-yy"
$ws.Range("H5").Value = 0.52083333333333337
$ws.Range("H5").NumberFormat = "h:mm"

$ws.Range("F5").Borders.Item(10).LineStyle = 1
$ws.Range("L5").Borders.Item(10).LineStyle = 1
$ws.Range("P5").Borders.Item(10).LineStyle = 1

# ---------------------------------------------------------------------------
# 9. Row 6 -- subject 4 (Olivia Hedge).
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = " Olivia Hedge "
$ws.Range("C6").Value = 16
$ws.Range("D6").Value = "Female"
$ws.Range("E6").Value = "olivia.hedge@gmail.com"
$ws.Range("G6").Value = 45315
$ws.Range("G6").NumberFormat = "mm-dd-yy"
$ws.Range("H6").Value = 0.10416666666666667
$ws.Range("H6").NumberFormat = "h:mm"

$ws.Range("F6").Borders.Item(10).LineStyle = 1
$ws.Range("L6").Borders.Item(10).LineStyle = 1
$ws.Range("P6").Borders.Item(10).LineStyle = 1

# ---------------------------------------------------------------------------
# 10. Selection / active cell to match the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("I7").Select()
